# Auto-generated edit script applying the Rafflesia_Profits market-data refresh
# (scheduled runner update) to the relevant job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 238.38095
$ws.Range("J33").Value = 752
$ws.Range("L33").Value = 752
$ws.Range("N33").Value = -1210
$ws.Range("H51").Value = 9989.9
$ws.Range("I51").Value = 9949.5
$ws.Range("K51").Value = 9949.5
$ws.Range("M51").Value = -9465.5
$ws.Range("H62").Value = 10998.8
$ws.Range("I62").Value = 9998
$ws.Range("K62").Value = 9998
$ws.Range("M62").Value = -9374
$ws.Range("H65").Value = 10998.8
$ws.Range("I65").Value = 9998
$ws.Range("K65").Value = 49990
$ws.Range("M65").Value = -46870
$ws.Range("H98").Value = 4998.5
$ws.Range("I98").Value = 4998.5
$ws.Range("K98").Value = 4998.5
$ws.Range("M98").Value = -3500.5
$ws.Range("H111").Value = 778.4286
$ws.Range("I111").Value = 775
$ws.Range("J111").Value = 783
$ws.Range("K111").Value = 2325
$ws.Range("L111").Value = 2349
$ws.Range("M111").Value = 742
$ws.Range("N111").Value = -8483
$ws.Range("H122").Value = 4998.5
$ws.Range("I122").Value = 4998.5
$ws.Range("K122").Value = 14995.5
$ws.Range("M122").Value = -12545.5
$ws.Range("H132").Value = 8844.333000000001
$ws.Range("I132").Value = 9371.286
$ws.Range("K132").Value = 28113.858
$ws.Range("M132").Value = -25583.858
$ws.Range("H137").Value = 2242.3809
$ws.Range("I137").Value = 2332.7778
$ws.Range("K137").Value = 6998.3334
$ws.Range("M137").Value = -4448.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9479.9
$ws.Range("I32").Value = 10488.777
$ws.Range("J32").Value = 400
$ws.Range("K32").Value = 10488.777
$ws.Range("L32").Value = 400
$ws.Range("M32").Value = -10201.777
$ws.Range("N32").Value = -974
$ws.Range("H43").Value = 70000
$ws.Range("J43").Value = 70000
$ws.Range("L43").Value = 70000
$ws.Range("N43").Value = -70626
$ws.Range("H61").Value = 3201
$ws.Range("I61").Value = 3201
$ws.Range("K61").Value = 3201
$ws.Range("M61").Value = -2989
$ws.Range("H74").Value = 4425
$ws.Range("I74").Value = 2991.6667
$ws.Range("K74").Value = 2991.6667
$ws.Range("M74").Value = -2117.6667
$ws.Range("H77").Value = 4425
$ws.Range("I77").Value = 2991.6667
$ws.Range("K77").Value = 14958.3335
$ws.Range("M77").Value = -10590.3335
$ws.Range("H122").Value = 1982.5454
$ws.Range("I122").Value = 2134.2222
$ws.Range("K122").Value = 6402.6666
$ws.Range("M122").Value = -3952.6666
$ws.Range("H127").Value = 17000
$ws.Range("I127").Value = 17000
$ws.Range("K127").Value = 17000
$ws.Range("M127").Value = -12040
$ws.Range("H132").Value = 5842.857
$ws.Range("I132").Value = 4375
$ws.Range("K132").Value = 13125
$ws.Range("M132").Value = -10595
$ws.Range("H136").Value = 3201
$ws.Range("I136").Value = 3201
$ws.Range("K136").Value = 9603
$ws.Range("M136").Value = -7053

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 22599.777
$ws.Range("I7").Value = 25174.75
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 25174.75
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -25061.75
$ws.Range("N7").Value = -2226
$ws.Range("H20").Value = 3965.2222
$ws.Range("I20").Value = 3965.2222
$ws.Range("K20").Value = 3965.2222
$ws.Range("M20").Value = -3718.2222
$ws.Range("H86").Value = 4246.0586
$ws.Range("I86").Value = 4309
$ws.Range("J86").Value = 3952.3333
$ws.Range("K86").Value = 4309
$ws.Range("L86").Value = 3952.3333
$ws.Range("M86").Value = -3186
$ws.Range("N86").Value = -6198.3333
$ws.Range("H89").Value = 4246.0586
$ws.Range("I89").Value = 4309
$ws.Range("J89").Value = 3952.3333
$ws.Range("K89").Value = 21545
$ws.Range("L89").Value = 19761.6665
$ws.Range("M89").Value = -15929
$ws.Range("N89").Value = -30993.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 34999
$ws.Range("J54").Value = 34999
$ws.Range("L54").Value = 34999
$ws.Range("N54").Value = -36315
$ws.Range("H88").Value = 6666.3335
$ws.Range("J88").Value = 6666.3335
$ws.Range("L88").Value = 6666.3335
$ws.Range("N88").Value = -7478.3335
$ws.Range("H91").Value = 6666.3335
$ws.Range("J91").Value = 6666.3335
$ws.Range("L91").Value = 6666.3335
$ws.Range("N91").Value = -9474.333500000001
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 1326
$ws.Range("I38").Value = 2540.6667
$ws.Range("J38").Value = 111.333336
$ws.Range("K38").Value = 7622.000100000001
$ws.Range("L38").Value = 334.000008
$ws.Range("M38").Value = -7275.000100000001
$ws.Range("N38").Value = -1028.000008
$ws.Range("H56").Value = 9999
$ws.Range("I56").Value = 9999
$ws.Range("K56").Value = 9999
$ws.Range("M56").Value = -9469
$ws.Range("H107").Value = 2751.8333
$ws.Range("J107").Value = 3001
$ws.Range("L107").Value = 9003
$ws.Range("N107").Value = -12843
$ws.Range("H113").Value = 1158
$ws.Range("I113").Value = 566
$ws.Range("K113").Value = 1698
$ws.Range("M113").Value = 472
$ws.Range("H122").Value = 1250
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5000
$ws.Range("I5").Value = 5000
$ws.Range("K5").Value = 5000
$ws.Range("M5").Value = -4888
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 11950
$ws.Range("J80").Value = 14891.667
$ws.Range("L80").Value = 14891.667
$ws.Range("N80").Value = -16887.667
$ws.Range("H83").Value = 11950
$ws.Range("J83").Value = 14891.667
$ws.Range("L83").Value = 74458.33499999999
$ws.Range("N83").Value = -84442.33499999999
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H132").Value = 1658.1666
$ws.Range("I132").Value = 1499.75
$ws.Range("J132").Value = 1975
$ws.Range("K132").Value = 4499.25
$ws.Range("L132").Value = 5925
$ws.Range("M132").Value = -1969.25
$ws.Range("N132").Value = -10985

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 700
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 800
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 800
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -1176
$ws.Range("H55").Value = 828.4545000000001
$ws.Range("I55").Value = 363.33334
$ws.Range("K55").Value = 363.33334
$ws.Range("M55").Value = -190.33334
$ws.Range("H93").Value = 862.25
$ws.Range("J93").Value = 450
$ws.Range("L93").Value = 450
$ws.Range("N93").Value = -2946
$ws.Range("H132").Value = 2499
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 7497
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -12557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3201
$ws.Range("I62").Value = 3002
$ws.Range("J62").Value = 3400
$ws.Range("K62").Value = 3002
$ws.Range("L62").Value = 3400
$ws.Range("M62").Value = -2378
$ws.Range("N62").Value = -4648
$ws.Range("H65").Value = 3201
$ws.Range("I65").Value = 3002
$ws.Range("J65").Value = 3400
$ws.Range("K65").Value = 15010
$ws.Range("L65").Value = 17000
$ws.Range("M65").Value = -11890
$ws.Range("N65").Value = -23240
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H113").Value = 1502.3
$ws.Range("I113").Value = 1017.5714
$ws.Range("J113").Value = 2633.3333
$ws.Range("K113").Value = 3052.7142
$ws.Range("L113").Value = 7899.999899999999
$ws.Range("M113").Value = -882.7142000000003
$ws.Range("N113").Value = -12239.9999
$ws.Range("H136").Value = 4124.4
$ws.Range("I136").Value = 4124.4
$ws.Range("K136").Value = 12373.2
$ws.Range("M136").Value = -9823.199999999999
